$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The account-statement table originally listed 2 workers x 2 arrears
# periods (1601 and 1602). This update removes the "1601" period rows,
# leaving only the "1602" period for each worker (part 1 of the new
# statement). Concretely that means dropping the old row 17 (WILLIAM /
# 1601) and the old row 18 (MILTON / 1602, duplicate-period slot), which
# shifts the former row 19 (MILTON / 1601, bottom-border styling) up to
# become the new row 17 - then that row's period value is corrected to
# 1602.
$ws.Rows(18).EntireRow.Delete()
$ws.Rows(17).EntireRow.Delete()

# Row 17 (was row 19) still shows the old "1601" period label - fix it to
# match the remaining data (1602).
$ws.Range("E17").Value = "1602"

# Update the totals to reflect the reduced data set.
$ws.Range("E11").Value = 51520
$ws.Range("F13").Value = 1
